# Auto-generated: updates cryptos list cell values to match the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.342.52"
$ws.Range("E2").Value = "  -4.77%  "
$ws.Range("D3").Value = "2.183.51"
$ws.Range("E3").Value = "  -7.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "480.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.91%  "
$ws.Range("D9").Value = "2.196.56"
$ws.Range("E9").Value = "  -7.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0910"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.312"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.56"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.36%  "
$ws.Range("D14").Value = "2.570.15"
$ws.Range("E14").Value = "  -7.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.88"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "53.252.51"
$ws.Range("E16").Value = "  -4.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000127"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.97%  "
$ws.Range("D18").Value = "2.171.50"
$ws.Range("E18").Value = "  -9.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "291.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.362"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "2.284.38"
$ws.Range("E27").Value = "  -7.75%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -7.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.64"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.811"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.51"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.363"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0872"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("E48").Value = "  -6.04%  "
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "226.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0198"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.04%  "
